# The ExcelDataSourceFile data has moved: the data that used to live in
# columns A:C (rows 1-3) now lives in columns F:H (rows 1-3). Move the
# cell contents accordingly (preserving the shared-string text), clear
# the old range, and update the active selection to the new location.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 1; $r -le 3; $r++) {
    for ($c = 1; $c -le 3; $c++) {
        $srcCell = $ws.Cells.Item($r, $c)
        $dstCell = $ws.Cells.Item($r, $c + 5)
        $dstCell.Value = $srcCell.Text
    }
}

# Remove the now-empty original range.
$ws.Range("A1:C3").Clear()

# Match the author's new selection over the relocated data.
$ws.Range("F1:H3").Select()
